$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.003.01'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.548.07'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.56'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.01'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.584'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -1.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.546.32'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -2.24%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.19'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +2.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.001.43'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +3.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.954.58'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.540.05'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.46'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.52'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.87'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('E25').Value = '  -3.39%  '
$ws.Range('E26').Value = '  +6.19%  '
$ws.Range('E27').Value = '  +11.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.38'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +3.00%  '
$ws.Range('E30').Value = '  +5.58%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0809'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.86'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.65'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E34').Value = '  +4.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '407.93'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +9.61%  '
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.37'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.07'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.24'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.76'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.604'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0957'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0518'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('E49').Value = '  +4.34%  '
$ws.Range('E50').Value = '  +1.09%  '
$ws.Range('E51').Value = '  +0.02%  '
